$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Text replacements (Key Skills table + Experience bullet)
#    Order matters: replace the old "TCP/IP networking fundamentals"
#    bullet before introducing the new one via the "AWS" replacement.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Linux (Debian)", $true, $true, $false, $false, $false, $true, 1, $false, "Remote desktop support", 2) | Out-Null
$d.Content.Find.Execute("Python", $true, $true, $false, $false, $false, $true, 1, $false, "MS Office Suite", 2) | Out-Null
$d.Content.Find.Execute("Bash", $true, $true, $false, $false, $false, $true, 1, $false, "Managing PHI in HIPAA regulated environment", 2) | Out-Null
$d.Content.Find.Execute("TCP/IP networking fundamentals", $true, $true, $false, $false, $false, $true, 1, $false, "User training", 2) | Out-Null
$d.Content.Find.Execute("AWS", $true, $true, $false, $false, $false, $true, 1, $false, "TCP/IP networking fundamentals", 2) | Out-Null
$d.Content.Find.Execute("Managed receipt and processing of sensitive medical data required for state/federal compliance audits", $true, $true, $false, $false, $false, $true, 1, $false, "Managed PHI receipt and processing from external medical organizations for HEDIS and QARR compliance", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Key Skills table geometry (2nd table in the document)
# ------------------------------------------------------------------
$t = $d.Tables.Item(2)

# Overall table width: 10080 dxa -> 10500 dxa (PreferredWidth is in points = dxa/20)
$t.PreferredWidth = 525

# Table indent: 100% -> 0%
$t2 = $d.Tables.Item(2)
$t2.Rows.LeftIndent = 0

# Column widths (dxa -> points = dxa/20): 4980/1515/3585 -> 3630/3795/3075
$t3 = $d.Tables.Item(2)
$t3.Columns.Item(1).Width = 181.5

$t4 = $d.Tables.Item(2)
$t4.Columns.Item(2).Width = 189.75

$t5 = $d.Tables.Item(2)
$t5.Columns.Item(3).Width = 153.75

# Cell margins: 100 dxa -> 0 dxa on all sides, for every cell in the row
$t6 = $d.Tables.Item(2)
for ($c = 1; $c -le $t6.Columns.Count; $c++) {
    $cell = $t6.Cell(1, $c)
    $cell.TopPadding = 0
    $cell.BottomPadding = 0
    $cell.LeftPadding = 0
    $cell.RightPadding = 0
}

Write-Output "done"
